# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45505 (2024-08-01) to 45506 (2024-08-02).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45505) {
        $cell.Value2 = 45506
    }
}
